$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row definitions matching the R script's appended output rows.
$rows = @(
    @{ Row = 180; Date = 45455.2916666667 },
    @{ Row = 181; Date = 45456.2916666667 }
)

# Scratch cell (far below the real data) used to mint a "text" cell whose
# value we can then copy by-value into column G without dragging a new
# number-format style onto the target cell itself.
$scratchRow = 1000
$ws.Cells.Item($scratchRow, 1).Value = "'4.23999977111816"
$ws.Cells.Item($scratchRow, 1).Copy()

foreach ($r in $rows) {
    $row = $r.Row

    # Clone row 179's formatting (the yyyy-mm-dd hh:mm:ss date style lives
    # on column A there) down onto the new row before writing values, so
    # the date cell reuses the existing style instead of minting a new one.
    $ws.Range("A179").Copy()
    $ws.Range("A$row").PasteSpecial(-4122)

    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 3).Value = 4.23999977111816
    $ws.Cells.Item($row, 4).Value = 4.23999977111816
    $ws.Cells.Item($row, 5).Value = 4.23999977111816
    $ws.Cells.Item($row, 6).Value = 4.23999977111816

    # Column G holds the adj_close value stored as text (matches the
    # existing rows 2-179, which are all shared-string "t=s" cells even
    # though they look numeric) - paste-by-value from the text scratch
    # cell so it lands as a shared string, not an auto-converted number.
    $ws.Cells.Item($scratchRow, 1).Copy()
    $ws.Cells.Item($row, 7).PasteSpecial(-4163)

    $ws.Cells.Item($row, 8).Value = "ELSA.MI"
}

$ws.Rows.Item($scratchRow).Delete()
